# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Alcachofa" (Española, Vega Central
# Mapocho de Santiago) at row 334, pushing the existing rows 334-353 down to
# 336-355.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 334 (shifts old rows 334:353 -> 336:355)
$ws.Range("A334:A335").EntireRow.Insert()

# --- New row 334: Española / Primera ---
$ws.Cells.Item(334, 1).Value = 9
$ws.Cells.Item(334, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(334, 3).Value = "Metropolitana"
$ws.Cells.Item(334, 4).Value = 44516
$ws.Cells.Item(334, 5).Value = 13
$ws.Cells.Item(334, 6).Value = 100112013
$ws.Cells.Item(334, 7).Value = "Alcachofa"
$ws.Cells.Item(334, 8).Value = "Española"
$ws.Cells.Item(334, 9).Value = "Primera"
$ws.Cells.Item(334, 10).Value = 52
$ws.Cells.Item(334, 11).Value = 9000
$ws.Cells.Item(334, 12).Value = 10000
$ws.Cells.Item(334, 13).Value = 9500
$ws.Cells.Item(334, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(334, 15).Value = "Región Metropolitana"
$ws.Cells.Item(334, 16).Value = 317
$ws.Cells.Item(334, 17).Value = 30
$ws.Cells.Item(334, 18).Value = "Hortaliza"

# --- New row 335: Española / Segunda ---
$ws.Cells.Item(335, 1).Value = 9
$ws.Cells.Item(335, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(335, 3).Value = "Metropolitana"
$ws.Cells.Item(335, 4).Value = 44516
$ws.Cells.Item(335, 5).Value = 13
$ws.Cells.Item(335, 6).Value = 100112013
$ws.Cells.Item(335, 7).Value = "Alcachofa"
$ws.Cells.Item(335, 8).Value = "Española"
$ws.Cells.Item(335, 9).Value = "Segunda"
$ws.Cells.Item(335, 10).Value = 34
$ws.Cells.Item(335, 11).Value = 7000
$ws.Cells.Item(335, 12).Value = 8000
$ws.Cells.Item(335, 13).Value = 7500
$ws.Cells.Item(335, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(335, 15).Value = "Región Metropolitana"
$ws.Cells.Item(335, 16).Value = 188
$ws.Cells.Item(335, 17).Value = 40
$ws.Cells.Item(335, 18).Value = "Hortaliza"

# Make sure the Date column keeps its date number format (style already
# inherited from the row above via Insert, but set it explicitly to be safe).
$ws.Range("D334:D335").NumberFormat = "YYYY-MM-DD HH:MM:SS"
